# Added a reset macro to the rubric
# (codeName / VBA project bits are session artifacts this host does not
# persist into the OOXML; we reproduce the observable content + formatting
# changes that the commit made to the workbook.)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Rubric"
$ws2 = $wb.Worksheets.Item(2)   # "Grade"

# --- sheet2 ("Grade"): rotate the three header/banner strings ----------
# Before: A1="Lab 2: A Review web site" (bold title)
#         A2 (merged A2:E2, tall wrap row) = "Excellent work!"
#         A3 (merged A3:E3, short row)     = "Here's the grade breakdown:"
# After:  A1="Here's the grade breakdown:"
#         A2 (merged A2:E2, tall wrap row) = "Lab 2: A Review web site"
#         A3 (merged A3:E3, short row)     = "Excellent work!"
$ws2.Range("A1").Value = "Here's the grade breakdown:"
$ws2.Range("A2").Value = "Lab 2: A Review web site"
$ws2.Range("A3").Value = "Excellent work!"

# Row 2's custom height shrinks a bit (38.4 -> 31.2)
$ws2.Rows.Item(2).RowHeight = 31.2

# Swap the "Possible"/"Actual" column headers on both mini-tables
$ws2.Range("B5").Value = "Possible"
$ws2.Range("C5").Value = "Actual"
$ws2.Range("B17").Value = "Possible"
$ws2.Range("C17").Value = "Actual"

# Column width tweaks
$ws2.Columns.Item(2).ColumnWidth = 7.796875
$ws2.Columns.Item(4).ColumnWidth = 0.69921875

# Selection moves from F9 to G8 on the Grade sheet
$ws2.Activate()
$ws2.Range("G8").Select()

# sheet1 ("Rubric") keeps its C7 selection; only the remembered scroll
# position (topLeftCell) is cleared, which happens on a normal re-save.
$ws1.Activate()
$ws1.Range("A1").Select()
$ws1.Range("C7").Select()

# Restore "Grade" as the active sheet/tab, matching the saved workbook view.
$ws2.Activate()
